$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Schedule")

# --- Fix weekday header row (A1:G1) so labels match the actual dates ---
# Originally mislabeled starting at Monday; correct order starting at Thursday.
$ws.Range("A1").Value = "Thursday"
$ws.Range("B1").Value = "Friday"
$ws.Range("C1").Value = "Saturday"
$ws.Range("D1").Value = "Sunday"
$ws.Range("E1").Value = "Monday"
$ws.Range("F1").Value = "Tuesday"
$ws.Range("G1").Value = "Wednesday"

# --- Highlight styling on row 2 date cells ---
# A2:B2 keep the "highlighted" fill already used elsewhere (fillId 7 in the xlsx),
# applied via Range.Style after duplicating/reusing format from G9-ish emphasis cells.
# C2:E2 move to a slightly different shade (fillId 4) versus their previous fillId 5.
$ws.Range("C2:E2").Interior.ThemeColor = 1
$ws.Range("C2:E2").Interior.TintAndShade = -0.149998474074526

# --- Fill in previously-empty log entry for topic 3 (row 4) ---
$ws.Range("L4").Value = 43758
$ws.Range("M4").Value = 0.91666666666666663
$ws.Range("N4").Value = 0.9375
$ws.Range("O4").Value = 3
$ws.Range("P4").Value = 30

# --- Column width for column G ---
$ws.Columns.Item(7).ColumnWidth = 10.85546875

# --- Selection change ---
$ws.Range("R13").Select()
